# Update model coefficient values (Wvpweerstand_modelcoefficienten)
# Each worksheet holds fitted model coefficients for a given signal;
# this refreshes offset/slope/temp/time/model_std/gewijzigd after re-fitting.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("IK106")
$ws.Range("B2").Value = -0.01436667179288742
$ws.Range("B4").Value = -0.0000008046470117331002
$ws.Range("B5").Value = 12.4336236677692
$ws.Range("B6").Value = 7.109866935233991
$ws.Range("B7").Value = 161.034200441859
$ws.Range("B10").Value = 0.296825660738558
$ws.Range("B11").Value = 45699.67276428798

$ws = $wb.Worksheets.Item("Q100")
$ws.Range("B2").Value = -0.01110803453187978
$ws.Range("B4").Value = -0.0000003983342205357038
$ws.Range("B5").Value = 11.61306265065054
$ws.Range("B6").Value = 5.258588144304044
$ws.Range("B7").Value = 160.4912958158204
$ws.Range("B10").Value = 0.2398145154032846
$ws.Range("B11").Value = 45699.67143572916

$ws = $wb.Worksheets.Item("Q200")
$ws.Range("B2").Value = -0.01106984620998736
$ws.Range("B4").Value = -0.0000001879784379687332
$ws.Range("B5").Value = 12.39714713199219
$ws.Range("B6").Value = 4.943663196731656
$ws.Range("B7").Value = 177.6318041255628
$ws.Range("B10").Value = 0.3475518427219827
$ws.Range("B11").Value = 45699.67150008102

$ws = $wb.Worksheets.Item("Q300")
$ws.Range("B2").Value = -0.01622926636488082
$ws.Range("B4").Value = -0.0000002638228465913171
$ws.Range("B5").Value = 12.16546815094099
$ws.Range("B6").Value = 6.680023662096612
$ws.Range("B7").Value = 173.6923360010707
$ws.Range("B10").Value = 0.2997107102016596
$ws.Range("B11").Value = 45699.67156946759

$ws = $wb.Worksheets.Item("Q400")
$ws.Range("B2").Value = -0.01248680672979542
$ws.Range("B4").Value = -0.0000005954965802022144
$ws.Range("B5").Value = 11.12271950505708
$ws.Range("B6").Value = 5.837740357183921
$ws.Range("B7").Value = 173.2577859748561
$ws.Range("B10").Value = 0.4171146902939639
$ws.Range("B11").Value = 45699.6716484375

$ws = $wb.Worksheets.Item("Q500")
$ws.Range("B2").Value = -0.01402050783430401
$ws.Range("B4").Value = -0.0000003506328900881021
$ws.Range("B5").Value = 10.76685570162649
$ws.Range("B6").Value = 4.116564695636098
$ws.Range("B7").Value = 180.5299384022
$ws.Range("B10").Value = 0.4416526914467693
$ws.Range("B11").Value = 45699.67173180555

$ws = $wb.Worksheets.Item("Q600")
$ws.Range("B2").Value = -0.01784971263347368
$ws.Range("B4").Value = -0.0000003572199002878696
$ws.Range("B5").Value = 11.26171424391456
$ws.Range("B6").Value = 3.431763141014778
$ws.Range("B7").Value = 170.135032131336
$ws.Range("B10").Value = 0.8847618532295183
$ws.Range("B11").Value = 45699.67179962963

$ws = $wb.Worksheets.Item("P100")
$ws.Range("B2").Value = -0.008807360459548873
$ws.Range("B4").Value = -0.00000008824247661175542
$ws.Range("B5").Value = 10.31238609600846
$ws.Range("B6").Value = 3.967781928581281
$ws.Range("B7").Value = 177.4547622180878
$ws.Range("B10").Value = 0.3730107261709022
$ws.Range("B11").Value = 45699.67187059027

$ws = $wb.Worksheets.Item("P200")
$ws.Range("B2").Value = -0.01279396309175958
$ws.Range("B4").Value = -0.0000006780263453958441
$ws.Range("B5").Value = 11.54835232707443
$ws.Range("B6").Value = 5.488265241023941
$ws.Range("B7").Value = 177.3975175370383
$ws.Range("B10").Value = 0.2576582829884954
$ws.Range("B11").Value = 45699.67196149306

$ws = $wb.Worksheets.Item("P300")
$ws.Range("B2").Value = -0.0120889506420294
$ws.Range("B4").Value = -0.0000000001000000000078997
$ws.Range("B5").Value = 9.66801273822297
$ws.Range("B6").Value = 3.87117962463416
$ws.Range("B7").Value = 156.4718267698503
$ws.Range("B10").Value = 1.006512531444407
$ws.Range("B11").Value = 45699.67202224537

$ws = $wb.Worksheets.Item("P400")
$ws.Range("B2").Value = -0.01606135392338998
$ws.Range("B4").Value = -0.00000007475878774398771
$ws.Range("B5").Value = 12.18231555679064
$ws.Range("B6").Value = 5.940654007156947
$ws.Range("B7").Value = 168.4933306795201
$ws.Range("B10").Value = 0.3534290661868898
$ws.Range("B11").Value = 45699.67208886574

$ws = $wb.Worksheets.Item("P500")
$ws.Range("B2").Value = -0.01696256126314312
$ws.Range("B4").Value = -0.000001235461859764849
$ws.Range("B5").Value = 11.68767004865758
$ws.Range("B6").Value = 6.51982840267387
$ws.Range("B7").Value = 164.8042277843232
$ws.Range("B10").Value = 0.261798062501053
$ws.Range("B11").Value = 45699.67217222222

$ws = $wb.Worksheets.Item("P600")
$ws.Range("B2").Value = -0.01274294607116233
$ws.Range("B4").Value = -0.000001186117429086021
$ws.Range("B5").Value = 9.638138853111498
$ws.Range("B6").Value = 1.719584105432685
$ws.Range("B7").Value = 161.1720462889251
$ws.Range("B10").Value = 0.679415172045599
$ws.Range("B11").Value = 45699.67224761574

$ws = $wb.Worksheets.Item("IK91")
$ws.Range("B2").Value = -0.03171207709568124
$ws.Range("B4").Value = -0.0000006640808991888091
$ws.Range("B5").Value = 13.12457404819691
$ws.Range("B6").Value = 4.850793108584173
$ws.Range("B7").Value = 192.7822173236272
$ws.Range("B10").Value = 0.4035749656714804
$ws.Range("B11").Value = 45699.67230025463

$ws = $wb.Worksheets.Item("IK92")
$ws.Range("B2").Value = -0.01217774837051229
$ws.Range("B4").Value = -0.0000007539840963555597
$ws.Range("B5").Value = 11.8920547605261
$ws.Range("B6").Value = 6.732531389095854
$ws.Range("B7").Value = 153.6475571887269
$ws.Range("B10").Value = 0.1330591400676031
$ws.Range("B11").Value = 45699.67233880787

$ws = $wb.Worksheets.Item("IK93")
$ws.Range("B2").Value = -0.0138499392906743
$ws.Range("B4").Value = -0.0000004625187281881448
$ws.Range("B5").Value = 13.10382884366287
$ws.Range("B6").Value = 7.395463308394059
$ws.Range("B7").Value = 156.2521507770053
$ws.Range("B10").Value = 0.187954995633264
$ws.Range("B11").Value = 45699.67237949074

$ws = $wb.Worksheets.Item("IK94")
$ws.Range("B2").Value = -0.01617768539418761
$ws.Range("B4").Value = -0.0000009606638622740289
$ws.Range("B5").Value = 12.07726498676792
$ws.Range("B6").Value = 7.101667672471773
$ws.Range("B7").Value = 157.9393208305148
$ws.Range("B10").Value = 0.2764838304688362
$ws.Range("B11").Value = 45699.67242075231

$ws = $wb.Worksheets.Item("IK95")
$ws.Range("B2").Value = -0.01931196065975306
$ws.Range("B4").Value = -0.000001616479033324338
$ws.Range("B5").Value = 13.13675760359978
$ws.Range("B6").Value = 5.950724724969729
$ws.Range("B7").Value = 147.3773540818518
$ws.Range("B10").Value = 0.6576064579900261
$ws.Range("B11").Value = 45699.6724678588

$ws = $wb.Worksheets.Item("IK96")
$ws.Range("B2").Value = -0.02005139699079833
$ws.Range("B4").Value = -0.0000032867996690011
$ws.Range("B5").Value = 13.67428404704193
$ws.Range("B6").Value = 4.915583932951407
$ws.Range("B7").Value = 208.8997151189449
$ws.Range("B10").Value = 0.2816166117686612
$ws.Range("B11").Value = 45699.67251216435

$ws = $wb.Worksheets.Item("IK101")
$ws.Range("B2").Value = -0.02500104128318027
$ws.Range("B4").Value = -0.0000005135898872162277
$ws.Range("B5").Value = 12.05143632470478
$ws.Range("B6").Value = 6.072419811520017
$ws.Range("B7").Value = 164.4690201006416
$ws.Range("B10").Value = 0.286387786735651
$ws.Range("B11").Value = 45699.67255212963

$ws = $wb.Worksheets.Item("IK102")
$ws.Range("B2").Value = -0.01309579322349816
$ws.Range("B5").Value = 12.43576456492437
$ws.Range("B6").Value = 6.177438809527719
$ws.Range("B7").Value = 153.2108751751385
$ws.Range("B10").Value = 0.19227165956517
$ws.Range("B11").Value = 45699.67259806713

$ws = $wb.Worksheets.Item("IK103")
$ws.Range("B2").Value = -0.01158812065812325
$ws.Range("B4").Value = -0.0000003938757397830367
$ws.Range("B5").Value = 12.29506440822717
$ws.Range("B6").Value = 5.031208770376971
$ws.Range("B7").Value = 144.8491858156742
$ws.Range("B10").Value = 0.2258624216083131
$ws.Range("B11").Value = 45699.67263927084

$ws = $wb.Worksheets.Item("IK104")
$ws.Range("B2").Value = -0.01385809094989918
$ws.Range("B4").Value = -0.0000003741895442759716
$ws.Range("B5").Value = 12.28859336398036
$ws.Range("B6").Value = 7.325898834286132
$ws.Range("B7").Value = 159.9939862223962
$ws.Range("B10").Value = 0.2219323447089247
$ws.Range("B11").Value = 45699.67267927084

$ws = $wb.Worksheets.Item("IK105")
$ws.Range("B2").Value = -0.01394070698381442
$ws.Range("B4").Value = -0.0000000001000000000000517
$ws.Range("B5").Value = 11.98984461152882
$ws.Range("B6").Value = 7.167725488565238
$ws.Range("B7").Value = 155.167126904195
$ws.Range("B10").Value = 0.2571575459472458
$ws.Range("B11").Value = 45699.67272045139
